$d = $word.ActiveDocument
$rng = $d.Range(416, 549)
$rng.Italic = 1
$rng.Text = "kregen behaalden betere resultaten. En van deze studenten behaalden de studenten die uitgesteld feedback kregen de beste resultaten."
$p1 = $d.Paragraphs.Item(1)

$newEnd = 416 + "kregen behaalden betere resultaten. En van deze studenten behaalden de studenten die uitgesteld feedback kregen de beste resultaten.".Length
Write-Output "newEnd=$newEnd  paraEnd=$($p1.Range.End)"
$fix = $d.Range(416, $newEnd)
$fix.Italic = 9999999
Write-Output $p1.Range.Text
